$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column before column D, shifting D:K -> E:L
$ws.Columns("D:D").Insert()

# Copy number formatting from column E (the old column D, now shifted) into new column D
# so the newly inserted cells carry the same style (date format row, number format rows, etc.)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2222100
$ws.Range("D9").Value = 1724700
$ws.Range("D10").Value = 497400
$ws.Range("D12").Value = 6300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 146000
$ws.Range("D15").Value = 224900
$ws.Range("D17").Value = 2353100
$ws.Range("D18").Value = -131000
$ws.Range("D20").Value = 2500
$ws.Range("D21").Value = 96300
$ws.Range("D22").Value = 3900
$ws.Range("D23").Value = -132400
$ws.Range("D24").Value = -2400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -130000
$ws.Range("D27").Value = -130000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -2500
$ws.Range("D33").Value = -130000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -130000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 135700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 309100
$ws.Range("D44").Value = 62600
$ws.Range("D45").Value = 22400
$ws.Range("D46").Value = 529800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 737300
$ws.Range("D49").Value = 115100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 42300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1424500
$ws.Range("D57").Value = 140100
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 104300
$ws.Range("D60").Value = 244400
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 26700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 271100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -120700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1153300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -130000
$ws.Range("D83").Value = 224900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 342100
$ws.Range("D91").Value = -311100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -276200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -44400
$ws.Range("D101").Value = 400
$ws.Range("D102").Value = 21900
